# Import vendedor (salesperson) data into the clientes workbook, with
# notification/contact columns O:R, including clickable mailto hyperlinks
# for the vendor e-mail column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O (cedula numbers) is stored as text so leading zeros are preserved.
# This has to be applied before any values are written into the column,
# otherwise Excel auto-converts numeric-looking text into real numbers.
$ws.Columns("O").NumberFormat = "@"

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("O1").Value = "vendedor_cedula"
$ws.Range("P1").Value = "vendedor_nombre"
$ws.Range("Q1").Value = "vendedor_apellido"
$ws.Range("R1").Value = "vendedor_email"

# --- Data rows --------------------------------------------------------------
# Entered in the same order the original author typed them (this keeps the
# shared-string table / cell ordering consistent with the source workbook).

# row 3 cedula entered first
$ws.Range("O3").Value = "09090900"

# row 2 (Juan Perez)
$ws.Range("P2").Value = "Juan"
$ws.Range("Q2").Value = "Perez"
$ws.Range("R2").Value = "perez@crm.com"

# row 3 (Claudio Loja)
$ws.Range("P3").Value = "Claudio"
$ws.Range("Q3").Value = "Loja"
$ws.Range("R3").Value = "juan.moscoso@primme.tech"

# row 2 cedula filled in afterwards
$ws.Range("O2").Value = "0103902399432"

# row 4 (nuevo)
$ws.Range("O4").Value = "0090902392"
$ws.Range("P4").Value = "algo"
$ws.Range("Q4").Value = "nuevo"
$ws.Range("R4").Value = "nuevo@otro.com"

# --- Hyperlinks for the vendor e-mail column ---------------------------------
$ws.Hyperlinks.Add($ws.Range("R2"), "mailto:perez@crm.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("R3"), "mailto:juan.moscoso@primme.tech") | Out-Null
$ws.Hyperlinks.Add($ws.Range("R4"), "mailto:nuevo@otro.com") | Out-Null

# --- View state: leave the active cell on O5 like in the source workbook ----
$ws.Range("O5").Select() | Out-Null
